# Generate Report for Handoff
#
# The a10bb7dc-... source file has now been handed off (status moves from
# "In Translation" to "Ready for handoff"), so its row swaps places with
# the already-handed-off 5e3adcf0-... row on every sheet, and its handoff
# file / handoff datetime are refreshed.

function Set-HyperlinkDisplay {
    param($ws, $rangeAddr, $newText)
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $rangeAddr) {
            $hl.TextToDisplay = $newText
        }
    }
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": rows for the two source files swap places, and the
# status for a10bb7dc... becomes "Ready for handoff" too.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "5e3adcf0-2dfd-4bd9-9041-3078e0acd852.md"
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
Set-HyperlinkDisplay $wsOverview '$A$2' "5e3adcf0-2dfd-4bd9-9041-3078e0acd852.md"

$wsOverview.Range("A3").Value = "a10bb7dc-2af9-4a10-ba96-0ca14e954cea.md"
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
Set-HyperlinkDisplay $wsOverview '$A$3' "a10bb7dc-2af9-4a10-ba96-0ca14e954cea.md"

# ---------------------------------------------------------------------
# Sheet "zh-cn": same row swap, plus refreshed handoff file / datetime
# for a10bb7dc... (now on row 3).
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "5e3adcf0-2dfd-4bd9-9041-3078e0acd852.md"
$wsZhCn.Range("B2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "5e3adcf0-2dfd-4bd9-9041-3078e0acd852.be8a2ac0bdfac6a0c26fa9cad11a1af283b2562e.zh-cn.xlf"
$wsZhCn.Range("D2").Value = "2016-02-22 13:45:48"
Set-HyperlinkDisplay $wsZhCn '$A$2' "5e3adcf0-2dfd-4bd9-9041-3078e0acd852.md"
Set-HyperlinkDisplay $wsZhCn '$C$2' "5e3adcf0-2dfd-4bd9-9041-3078e0acd852.be8a2ac0bdfac6a0c26fa9cad11a1af283b2562e.zh-cn.xlf"

$wsZhCn.Range("A3").Value = "a10bb7dc-2af9-4a10-ba96-0ca14e954cea.md"
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("C3").Value = "a10bb7dc-2af9-4a10-ba96-0ca14e954cea.2d3f83a0351f35ed02ad7ecee08b2de4891d8c6e.zh-cn.xlf"
$wsZhCn.Range("D3").Value = "2016-02-22 13:48:41"
Set-HyperlinkDisplay $wsZhCn '$A$3' "a10bb7dc-2af9-4a10-ba96-0ca14e954cea.md"
Set-HyperlinkDisplay $wsZhCn '$C$3' "a10bb7dc-2af9-4a10-ba96-0ca14e954cea.2d3f83a0351f35ed02ad7ecee08b2de4891d8c6e.zh-cn.xlf"

# ---------------------------------------------------------------------
# Sheet "de-de": same row swap, plus refreshed handoff file / datetime
# for a10bb7dc... (now on row 3).
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "5e3adcf0-2dfd-4bd9-9041-3078e0acd852.md"
$wsDeDe.Range("B2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "5e3adcf0-2dfd-4bd9-9041-3078e0acd852.be8a2ac0bdfac6a0c26fa9cad11a1af283b2562e.de-de.xlf"
$wsDeDe.Range("D2").Value = "2016-02-22 13:46:04"
Set-HyperlinkDisplay $wsDeDe '$A$2' "5e3adcf0-2dfd-4bd9-9041-3078e0acd852.md"
Set-HyperlinkDisplay $wsDeDe '$C$2' "5e3adcf0-2dfd-4bd9-9041-3078e0acd852.be8a2ac0bdfac6a0c26fa9cad11a1af283b2562e.de-de.xlf"

$wsDeDe.Range("A3").Value = "a10bb7dc-2af9-4a10-ba96-0ca14e954cea.md"
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("C3").Value = "a10bb7dc-2af9-4a10-ba96-0ca14e954cea.2d3f83a0351f35ed02ad7ecee08b2de4891d8c6e.de-de.xlf"
$wsDeDe.Range("D3").Value = "2016-02-22 13:48:57"
Set-HyperlinkDisplay $wsDeDe '$A$3' "a10bb7dc-2af9-4a10-ba96-0ca14e954cea.md"
Set-HyperlinkDisplay $wsDeDe '$C$3' "a10bb7dc-2af9-4a10-ba96-0ca14e954cea.2d3f83a0351f35ed02ad7ecee08b2de4891d8c6e.de-de.xlf"
